$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M8").Value = 1.05
$ws.Range("O8").Value = 1.29
$ws.Range("M9").Value = 1.04
$ws.Range("O9").Value = 1.22
$ws.Range("G13").Value = 2.15
$ws.Range("N13").Value = 17
$ws.Range("AE13").Value = 12
$ws.Range("AZ13").Value = 51
$ws.Range("BC13").Value = 351
$ws.Range("M15").Value = 1.05
$ws.Range("N15").Value = 11
$ws.Range("G17").Value = 2.75
$ws.Range("H17").Value = 3.25
$ws.Range("I17").Value = 2.38
$ws.Range("J17").Value = 3.5
$ws.Range("K17").Value = 2.1
$ws.Range("L17").Value = 3.1
$ws.Range("M17").Value = 1.06
$ws.Range("N17").Value = 10
$ws.Range("Q17").Value = 2.08
$ws.Range("R17").Value = 1.73
$ws.Range("W17").Value = 8.5
$ws.Range("Z17").Value = 29
$ws.Range("AC17").Value = 9.5
$ws.Range("AH17").Value = 11
$ws.Range("AJ17").Value = 23
$ws.Range("AL17").Value = 29
$ws.Range("AM17").Value = 251
$ws.Range("AO17").Value = 17
$ws.Range("AX17").Value = 13
$ws.Range("AY17").Value = 23
$ws.Range("AZ17").Value = 41
$ws.Range("BB17").Value = 151
$ws.Range("I23").Value = 2.8
$ws.Range("J23").Value = 3.25
$ws.Range("AK23").Value = 23
$ws.Range("AR23").Value = 81
$ws.Range("N24").Value = 9
$ws.Range("G28").Value = 2.52
$ws.Range("H28").Value = 3.15
$ws.Range("I28").Value = 2.57
$ws.Range("J28").Value = 3.2
$ws.Range("K28").Value = 2.07
$ws.Range("L28").Value = 3.2
$ws.Range("N28").Value = 7.2
$ws.Range("S28").Value = 1.42
$ws.Range("T28").Value = 2.67
$ws.Range("U28").Value = 1.7
$ws.Range("V28").Value = 2.05
$ws.Range("X28").Value = 13
$ws.Range("Z28").Value = 28
$ws.Range("AC28").Value = 7.2
$ws.Range("AD28").Value = 6.2
$ws.Range("AE28").Value = 13
$ws.Range("AF28").Value = 55
$ws.Range("AH28").Value = 13.5
$ws.Range("AJ28").Value = 30
$ws.Range("AL28").Value = 29
$ws.Range("AM28").Value = 400
$ws.Range("AO28").Value = 14
$ws.Range("AP28").Value = 22
$ws.Range("AR28").Value = 100
$ws.Range("AT28").Value = 2.67
$ws.Range("AU28").Value = 6.9
$ws.Range("AV28").Value = 60
$ws.Range("AZ28").Value = 65
$ws.Range("G30").Value = 1.47
$ws.Range("H30").Value = 4.45
$ws.Range("I30").Value = 5
$ws.Range("J30").Value = 1.91
$ws.Range("K30").Value = 2.62
$ws.Range("L30").Value = 4.8
$ws.Range("N30").Value = 10
$ws.Range("T30").Value = 3.8
$ws.Range("V30").Value = 2.32
$ws.Range("Y30").Value = 8.5
$ws.Range("Z30").Value = 12
$ws.Range("AB30").Value = 18
$ws.Range("AC30").Value = 10
$ws.Range("AD30").Value = 9.75
$ws.Range("AE30").Value = 14
$ws.Range("AF30").Value = 45
$ws.Range("AM30").Value = 250
$ws.Range("AN30").Value = 3.8
$ws.Range("AO30").Value = 6.8
$ws.Range("AQ30").Value = 17.5
$ws.Range("AT30").Value = 3.8
$ws.Range("AU30").Value = 7
$ws.Range("AV30").Value = 45
$ws.Range("AY30").Value = 24
$ws.Range("BA30").Value = 120
